# Management_information_DKI_2016.xlsx update
# - Standardize fertilizer naming: "KAS" -> "CAN", "Alzon 40 + 5 S" -> "UAN (inhibited)"
# - Fill in previously missing application dates for the 2nd/3rd nitrogen applications
# - Remove the erroneous "others" / "Bittersalz" fertilization rows that duplicated
#   data which belongs elsewhere (rows 50-52, columns B:E)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fertilizer name standardization -------------------------------------------------
$ws.Range("E47").Value = "CAN"
$ws.Range("E48").Value = "CAN"
$ws.Range("E49").Value = "UAN (inhibited)"

# --- Fill in missing application dates -------------------------------------------------
$ws.Range("C48").Value = 42473
$ws.Range("C49").Value = 42474

# --- Remove stray / duplicated fertilization entries -------------------------------------------------
$ws.Range("B50").Value = ""
$ws.Range("C50").Value = ""
$ws.Range("E50").Value = ""

$ws.Range("C51").Value = ""
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = ""

$ws.Range("C52").Value = ""
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = ""
